# Add the new "Angelo" user record to the Usuarios sheet.
# Columns: A=Nome, B=Usuario, C=Senha, D=CPF, E=Nascimento
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")
$ws.Activate()

$ws.Range("A1").Value = "Angelo"
$ws.Range("B1").Value = "Angelo"

# CPF/Senha/Nascimento look numeric/date-like to Excel's input parser, so a
# leading apostrophe forces them to stay plain text (matching the source
# data, which stores them as shared strings) - then drop the resulting
# "quote prefix" cell format so only the values/types change.
$ws.Range("C1").Value = "'123"
$ws.Range("D1").Value = "'96655682215"
$ws.Range("E1").Value = "'03/02/2005"
$ws.Range("C1:E1").ClearFormats()
